$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is plain text in the source workbook (t="inlineStr").
# Setting .Value directly lets Excel auto-detect some strings as numbers
# (dropping things like trailing zeros), so we force Text format, assign
# the literal string, then clear the temporary format back off again so
# the cell keeps its original (default/no-style) formatting.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '62.394.11'
Set-TextValue $ws.Range('E2') '  -2.49%  '
Set-TextValue $ws.Range('D3') '3.186.72'
Set-TextValue $ws.Range('E3') '  -3.90%  '
Set-TextValue $ws.Range('E4') '  -0.02%  '
Set-TextValue $ws.Range('D5') '587.18'
Set-TextValue $ws.Range('E5') '  -2.29%  '
Set-TextValue $ws.Range('D6') '135.07'
Set-TextValue $ws.Range('E6') '  -5.86%  '
Set-TextValue $ws.Range('E7') '  -0.06%  '
Set-TextValue $ws.Range('D8') '3.185.65'
Set-TextValue $ws.Range('E8') '  -3.95%  '
Set-TextValue $ws.Range('E9') '  -4.21%  '
Set-TextValue $ws.Range('D10') '0.142'
Set-TextValue $ws.Range('E10') '  -5.60%  '
Set-TextValue $ws.Range('D11') '5.25'
Set-TextValue $ws.Range('E11') '  -5.44%  '
Set-TextValue $ws.Range('D12') '0.450'
Set-TextValue $ws.Range('E12') '  -5.14%  '
Set-TextValue $ws.Range('D13') '0.0000235'
Set-TextValue $ws.Range('E13') '  -6.28%  '
Set-TextValue $ws.Range('D14') '33.21'
Set-TextValue $ws.Range('E14') '  -5.06%  '
Set-TextValue $ws.Range('D15') '3.708.72'
Set-TextValue $ws.Range('E15') '  -4.07%  '
Set-TextValue $ws.Range('E16') '  -1.26%  '
Set-TextValue $ws.Range('D17') '3.182.42'
Set-TextValue $ws.Range('E17') '  -4.05%  '
Set-TextValue $ws.Range('D18') '62.423.48'
Set-TextValue $ws.Range('E18') '  -2.57%  '
Set-TextValue $ws.Range('D19') '6.58'
Set-TextValue $ws.Range('E19') '  -4.93%  '
Set-TextValue $ws.Range('D20') '456.60'
Set-TextValue $ws.Range('E20') '  -5.45%  '
Set-TextValue $ws.Range('D21') '13.86'
Set-TextValue $ws.Range('E21') '  -3.37%  '
Set-TextValue $ws.Range('E22') '  -4.84%  '
Set-TextValue $ws.Range('D23') '7.61'
Set-TextValue $ws.Range('E23') '  -4.94%  '
Set-TextValue $ws.Range('D24') '13.41'
Set-TextValue $ws.Range('E24') '  -2.08%  '
Set-TextValue $ws.Range('D25') '82.51'
Set-TextValue $ws.Range('E25') '  -2.94%  '
Set-TextValue $ws.Range('B27') 'PancakeSwap'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D27') '2.68'
Set-TextValue $ws.Range('E27') '  -3.85%  '
Set-TextValue $ws.Range('B28') 'FirstDigitalUSD'
Set-TextValue $ws.Range('C28') 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D28') '1.00'
Set-TextValue $ws.Range('E28') '  -0.07%  '
Set-TextValue $ws.Range('D29') '6.90'
Set-TextValue $ws.Range('E29') '  -5.50%  '
Set-TextValue $ws.Range('D30') '7.81'
Set-TextValue $ws.Range('E30') '  -5.45%  '
Set-TextValue $ws.Range('E31') '  -7.36%  '
Set-TextValue $ws.Range('D32') '27.29'
Set-TextValue $ws.Range('E32') '  -8.19%  '
Set-TextValue $ws.Range('E33') '  -4.38%  '
Set-TextValue $ws.Range('E34') '  -7.50%  '
Set-TextValue $ws.Range('E35') '  -6.16%  '
Set-TextValue $ws.Range('D36') '5.80'
Set-TextValue $ws.Range('E36') '  -3.43%  '
Set-TextValue $ws.Range('D37') '51.19'
Set-TextValue $ws.Range('D38') '0.0₃0691'
Set-TextValue $ws.Range('E38') '  -9.16%  '
Set-TextValue $ws.Range('D39') '0.0386'
Set-TextValue $ws.Range('E39') '  -4.28%  '
Set-TextValue $ws.Range('D40') '410.79'
Set-TextValue $ws.Range('E40') '  -5.78%  '
Set-TextValue $ws.Range('D41') '2.940.12'
Set-TextValue $ws.Range('E41') '  -4.21%  '
Set-TextValue $ws.Range('E42') '  +1.13%  '
Set-TextValue $ws.Range('D43') '8.01'
Set-TextValue $ws.Range('E43') '  -5.32%  '
Set-TextValue $ws.Range('D44') '2.63'
Set-TextValue $ws.Range('E44') '  -5.91%  '
Set-TextValue $ws.Range('E45') '  -3.33%  '
Set-TextValue $ws.Range('D46') '0.250'
Set-TextValue $ws.Range('E46') '  -6.79%  '
Set-TextValue $ws.Range('D47') '35.99'
Set-TextValue $ws.Range('E47') '  -0.01%  '
Set-TextValue $ws.Range('D48') '0.998'
Set-TextValue $ws.Range('E48') '  -0.12%  '
Set-TextValue $ws.Range('D49') '25.50'
Set-TextValue $ws.Range('E49') '  -4.21%  '
Set-TextValue $ws.Range('D50') '124.00'
Set-TextValue $ws.Range('E50') '  +0.49%  '
Set-TextValue $ws.Range('E51') '  -4.39%  '
